# WorkerCalc.xlsx edit: tweak worker-capacity multipliers/bases, switch the
# F-column formula (rows 15:26) from a straight MIN() to a FLOOR-averaged
# MIN(), and extend the second table down to rows 25 and 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top table (rows 4:12): F4 multiplier 3 -> 1 (propagates via F5=F4, F6=F5, ...) ---
$ws.Range("F4").Value = 1

# --- Second table (rows 15:26) ---

# Row 15 base inputs
$ws.Range("E15").Value = 4
$ws.Range("G15").Value = 0.5

# F column formula changed on every row from 5*MIN(C,D,6) to 5*MIN(FLOOR((C+D)/2,1),6)
$ws.Range("F15").Formula = "=5*MIN(FLOOR((C15+D15)/2, 1),6)"
$ws.Range("F16").Formula = "=5*MIN(FLOOR((C16+D16)/2, 1),6)"
$ws.Range("F17").Formula = "=5*MIN(FLOOR((C17+D17)/2, 1),6)"
$ws.Range("F18").Formula = "=5*MIN(FLOOR((C18+D18)/2, 1),6)"
$ws.Range("F19").Formula = "=5*MIN(FLOOR((C19+D19)/2, 1),6)"
$ws.Range("F20").Formula = "=5*MIN(FLOOR((C20+D20)/2, 1),6)"
$ws.Range("F21").Formula = "=5*MIN(FLOOR((C21+D21)/2, 1),6)"
$ws.Range("F22").Formula = "=5*MIN(FLOOR((C22+D22)/2, 1),6)"
$ws.Range("F23").Formula = "=5*MIN(FLOOR((C23+D23)/2, 1),6)"
$ws.Range("F24").Formula = "=5*MIN(FLOOR((C24+D24)/2, 1),6)"

# Row 17 width input
$ws.Range("D17").Value = 4

# Row 19 / 20 width inputs swapped (5/4 -> 4/5)
$ws.Range("D19").Value = 4
$ws.Range("D20").Value = 5

# Row 24: household counts shrunk, and E24 now mirrors E22 instead of E23
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 9
$ws.Range("E24").Formula = "=E22"

# --- New rows 25 and 26 ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 10
$ws.Range("E25").Formula = "=E23"
$ws.Range("F25").Formula = "=5*MIN(FLOOR((C25+D25)/2, 1),6)"
$ws.Range("G25").Formula = "=G24"
$ws.Range("H25").Formula = "=(C25*extra)*(D25*extra)*E25*F25*G25"

$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 14
$ws.Range("E26").Formula = "=E25"
$ws.Range("F26").Formula = "=5*MIN(FLOOR((C26+D26)/2, 1),6)"
$ws.Range("G26").Formula = "=G25"
$ws.Range("H26").Formula = "=(C26*extra)*(D26*extra)*E26*F26*G26"

# Selection moves to D21 (last cell the author touched)
$ws.Range("D21").Select()
